$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.526.69"
$ws.Range("E2").Value = "  -5.46%  "
$ws.Range("D3").Value = "3.467.62"
$ws.Range("E3").Value = "  -7.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.96%  "
$ws.Range("D7").Value = "3.463.69"
$ws.Range("E7").Value = "  -7.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.31%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.646"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -11.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.140"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -13.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -14.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -14.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -11.09%  "
$ws.Range("D15").Value = "4.035.81"
$ws.Range("E15").Value = "  -6.71%  "
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").Value = "3.469.01"
$ws.Range("E17").Value = "  -7.07%  "
$ws.Range("D18").Value = "65.385.65"
$ws.Range("E18").Value = "  -5.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.08%  "
$ws.Range("E21").Value = "  -11.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.90%  "
$ws.Range("E26").Value = "  -9.14%  "
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.98%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.49%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "605.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "62.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.69%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.391"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "0.0₃0707"
$ws.Range("E41").Value = "  -16.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.128"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.33%  "
$ws.Range("D43").Value = "2.934.26"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("E44").Value = "  -11.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.40%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("E47").Value = "  -12.67%  "
$ws.Range("E48").Value = "  -9.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "136.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -12.30%  "
